# "Add second office hw" - append a new shipper/product and two new price
# rows to the tender workbook's three sheets.

$wb = $excel.ActiveWorkbook

# --- shipper: new row 4 (S03 / Max / 10 / max@com.ua) -----------------
$shipper = $wb.Worksheets.Item("shipper")
$shipper.Range("A4").Value = "S03"
$shipper.Range("B4").Value = "Max"
# Force text so "10" isn't auto-coerced to a number (matches the other
# rating cells, which are stored as text), then drop back to the default
# style so no stray number-format style is left behind.
$shipper.Range("C4").NumberFormat = "@"
$shipper.Range("C4").Value = "10"
$shipper.Range("C4").Style = "Normal"
$shipper.Range("D4").Value = "max@com.ua"

# --- product: new row 4 (P03 / Папір) ----------------------------------
$product = $wb.Worksheets.Item("product")
$product.Range("A4").Value = "P03"
$product.Range("B4").Value = "Папір"

# --- price: new rows 4 and 5 -------------------------------------------
$price = $wb.Worksheets.Item("price")

$price.Range("A4").Value = "S02"
$price.Range("B4").Value = "P01"
$price.Range("C4").Value = "2,6"
$price.Range("D4").NumberFormat = "@"
$price.Range("D4").Value = "5"
$price.Range("D4").Style = "Normal"

$price.Range("A5").Value = "S03"
$price.Range("B5").Value = "P01"
$price.Range("C5").Value = "3,0"
$price.Range("D5").NumberFormat = "@"
$price.Range("D5").Value = "6"
$price.Range("D5").Style = "Normal"
